$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price (column D) and 1h volume-change (column E) figures
# for the Fri Jan 13 22:47:17 UTC 2023 GitHub Actions symbol-list refresh.
#
# The source cells are plain inline-string "Text" cells (t="inlineStr"),
# with no explicit style. A bare numeric-looking assignment like
# $range.Value = "293.12" gets auto-coerced by Excel into a Number cell,
# which would change the cell's stored type/format and diverge from the
# diff. Prefixing the literal with a single quote forces Excel to keep it
# as Text (the standard "quote prefix" trick); resetting Style to
# "Normal" afterwards strips the quote-prefix formatting flag back off so
# the cell's style index is left exactly as it was before the edit.

$updates = [ordered]@{
    "D2" = "293.12"
    "E2" = "2.27%"
    "D3" = "29.54"
    "E3" = "3.23%"
    "D4" = "5.253"
    "E4" = "3.68%"
    "D5" = "0.07161"
    "E5" = "7.77%"
    "D6" = "7.537"
    "E6" = "2.32%"
    "D7" = "3.596"
    "E7" = "5.59%"
    "E8" = "2.57%"
    "D9" = "0.9115"
    "E9" = "-2.94%"
    "D10" = "0.1633"
    "E10" = "3.52%"
    "D11" = "0.07651"
    "E11" = "15.80%"
    "D12" = "0.07765"
    "E12" = "2.54%"
    "D13" = "0.02918"
    "E13" = "-0.82%"
    "D14" = "0.08997"
    "E14" = "0.25%"
    "D15" = "0.001593"
    "E15" = "-0.40%"
    "E16" = "1.58%"
    "D17" = "0.006086"
    "E17" = "-3.00%"
    "D18" = "3.486"
    "E18" = "1.26%"
    "D20" = "0.3251"
    "E20" = "1.06%"
    "D21" = "0.1369"
    "E21" = "5.55%"
    "D22" = "4.047"
    "E22" = "-0.60%"
    "E23" = "2.51%"
    "D24" = "0.04516"
    "E24" = "0.49%"
    "D25" = "0.001206"
    "E25" = "2.07%"
    "D26" = "0.004268"
    "E26" = "2.92%"
    "E27" = "-6.60%"
    "D28" = "0.0001679"
    "E28" = "3.86%"
    "D40" = "0.04406"
    "E40" = "4.90%"
    "E41" = "4.21%"
    "D42" = "0.1281"
    "E42" = "2.45%"
    "E43" = "9.17%"
    "D44" = "0.01333"
    "E44" = "7.99%"
    "D45" = "0.00005825"
    "E45" = "4.24%"
    "D47" = "0.01292"
    "E47" = "-1.10%"
}

foreach ($cellRef in $updates.Keys) {
    $newText = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $newText
    $range.Style = "Normal"
}
